$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend header row (row 1) with two more sequence values ---
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
# Match formatting of the existing header cells (bold, bordered, centered)
$ws.Range("O1").Copy() | Out-Null
$ws.Range("P1:Q1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# --- For every data row (2-25), fix the I/K/M/O pattern and add P/Q columns ---
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value  = 2   # I -> 2
    $ws.Cells.Item($r, 11).Value = 1   # K -> 1
    $ws.Cells.Item($r, 13).Value = 2   # M -> 2
    $ws.Cells.Item($r, 15).Value = 1   # O -> 1
    $ws.Cells.Item($r, 16).Value = 2   # P (new)
    $ws.Cells.Item($r, 17).Value = 2   # Q (new)
}
